$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value2 = "NEW"
$ws.Range("I16").Value2 = "VACCATED"
$ws.Range("A40").Value2 = "VERONICA VACCATED"
$ws.Range("A38").Value2 = "PAID ON 15/12"
$ws.Range("E38").Value2 = "PAID ON 15/12"
